$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.320746
$ws.Range("H2").Value = 288.962238
$ws.Range("I2").Value = 0.3809824610908788
$ws.Range("J2").Value = 0.3809824610908788
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.935295
$ws.Range("N2").Value = 8.805885
$ws.Range("O2").Value = 0.0527739323334782
$ws.Range("P2").Value = 0.0527739323334782
$ws.Range("Q2").Value = 282.72980413007
$ws.Range("R2").Value = 2544.56823717063
$ws.Range("S2").Value = 0.02010594262185203
$ws.Range("T2").Value = 0.02010594262185203

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.320746
$ws.Range("H3").Value = 288.962238
$ws.Range("I3").Value = 0.3809824610908788
$ws.Range("J3").Value = 0.3809824610908788
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 38.46708166666667
$ws.Range("N3").Value = 115.401245
$ws.Range("O3").Value = 0.69160311482936
$ws.Range("P3").Value = 0.69160311482936
$ws.Range("Q3").Value = 3705.178002576256
$ws.Range("R3").Value = 33346.60202318631
$ws.Range("S3").Value = 0.2634886567858072
$ws.Range("T3").Value = 0.2634886567858072

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.320746
$ws.Range("H4").Value = 288.962238
$ws.Range("I4").Value = 0.3809824610908788
$ws.Range("J4").Value = 0.3809824610908788
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.981753333333333
$ws.Range("N4").Value = 17.94526
$ws.Range("O4").Value = 0.1075464802171131
$ws.Range("P4").Value = 0.1075464802171131
$ws.Range("Q4").Value = 576.1669434546533
$ws.Range("R4").Value = 5185.502491091879
$ws.Range("S4").Value = 0.04097332271477725
$ws.Range("T4").Value = 0.04097332271477725

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 96.320746
$ws.Range("H5").Value = 288.962238
$ws.Range("I5").Value = 0.3809824610908788
$ws.Range("J5").Value = 0.3809824610908788
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.236038333333333
$ws.Range("N5").Value = 24.708115
$ws.Range("O5").Value = 0.1480764726200487
$ws.Range("P5").Value = 0.1480764726200487
$ws.Range("Q5").Value = 793.3013563512633
$ws.Range("R5").Value = 7139.71220716137
$ws.Range("S5").Value = 0.05641453896844228
$ws.Range("T5").Value = 0.05641453896844229

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.46467533333333
$ws.Range("H6").Value = 55.394026
$ws.Range("I6").Value = 0.07303429161291354
$ws.Range("J6").Value = 0.07303429161291354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.935295
$ws.Range("N6").Value = 8.805885
$ws.Range("O6").Value = 0.0527739323334782
$ws.Range("P6").Value = 0.0527739323334782
$ws.Range("Q6").Value = 54.19926918255666
$ws.Range("R6").Value = 487.79342264301
$ws.Range("S6").Value = 0.003854306763603414
$ws.Range("T6").Value = 0.003854306763603414

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.46467533333333
$ws.Range("H7").Value = 55.394026
$ws.Range("I7").Value = 0.07303429161291354
$ws.Range("J7").Value = 0.07303429161291354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 38.46708166666667
$ws.Range("N7").Value = 115.401245
$ws.Range("O7").Value = 0.69160311482936
$ws.Range("P7").Value = 0.69160311482936
$ws.Range("Q7").Value = 710.2821739958188
$ws.Range("R7").Value = 6392.539565962369
$ws.Range("S7").Value = 0.05051074356884681
$ws.Range("T7").Value = 0.05051074356884681

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.46467533333333
$ws.Range("H8").Value = 55.394026
$ws.Range("I8").Value = 0.07303429161291354
$ws.Range("J8").Value = 0.07303429161291354
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.981753333333333
$ws.Range("N8").Value = 17.94526
$ws.Range("O8").Value = 0.1075464802171131
$ws.Range("P8").Value = 0.1075464802171131
$ws.Range("Q8").Value = 110.4511332240844
$ws.Range("R8").Value = 994.0601990167598
$ws.Range("S8").Value = 0.007854580998119073
$ws.Range("T8").Value = 0.007854580998119073

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.46467533333333
$ws.Range("H9").Value = 55.394026
$ws.Range("I9").Value = 0.07303429161291354
$ws.Range("J9").Value = 0.07303429161291354
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.236038333333333
$ws.Range("N9").Value = 24.708115
$ws.Range("O9").Value = 0.1480764726200487
$ws.Range("P9").Value = 0.1480764726200487
$ws.Range("Q9").Value = 152.0757738578878
$ws.Range("R9").Value = 1368.68196472099
$ws.Range("S9").Value = 0.01081466028234424
$ws.Range("T9").Value = 0.01081466028234425

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 122.909391
$ws.Range("H10").Value = 368.728173
$ws.Range("I10").Value = 0.4861499128584522
$ws.Range("J10").Value = 0.4861499128584522
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.935295
$ws.Range("N10").Value = 8.805885
$ws.Range("O10").Value = 0.0527739323334782
$ws.Range("P10").Value = 0.0527739323334782
$ws.Range("Q10").Value = 360.7753208553449
$ws.Range("R10").Value = 3246.977887698104
$ws.Range("S10").Value = 0.02565604260511828
$ws.Range("T10").Value = 0.02565604260511828

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 122.909391
$ws.Range("H11").Value = 368.728173
$ws.Range("I11").Value = 0.4861499128584522
$ws.Range("J11").Value = 0.4861499128584522
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 38.46708166666667
$ws.Range("N11").Value = 115.401245
$ws.Range("O11").Value = 0.69160311482936
$ws.Range("P11").Value = 0.69160311482936
$ws.Range("Q11").Value = 4727.965581197264
$ws.Range("R11").Value = 42551.69023077538
$ws.Range("S11").Value = 0.3362227940069275
$ws.Range("T11").Value = 0.3362227940069275

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 122.909391
$ws.Range("H12").Value = 368.728173
$ws.Range("I12").Value = 0.4861499128584522
$ws.Range("J12").Value = 0.4861499128584522
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 5.981753333333333
$ws.Range("N12").Value = 17.94526
$ws.Range("O12").Value = 0.1075464802171131
$ws.Range("P12").Value = 0.1075464802171131
$ws.Range("Q12").Value = 735.2136593122199
$ws.Range("R12").Value = 6616.922933809979
$ws.Range("S12").Value = 0.05228371198578277
$ws.Range("T12").Value = 0.05228371198578277

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 122.909391
$ws.Range("H13").Value = 368.728173
$ws.Range("I13").Value = 0.4861499128584522
$ws.Range("J13").Value = 0.4861499128584522
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.236038333333333
$ws.Range("N13").Value = 24.708115
$ws.Range("O13").Value = 0.1480764726200487
$ws.Range("P13").Value = 0.1480764726200487
$ws.Range("Q13").Value = 1012.286455802655
$ws.Range("R13").Value = 9110.578102223893
$ws.Range("S13").Value = 0.07198736426062365
$ws.Range("T13").Value = 0.07198736426062366

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.127183
$ws.Range("H14").Value = 45.381549
$ws.Range("I14").Value = 0.05983333443775553
$ws.Range("J14").Value = 0.05983333443775553
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.935295
$ws.Range("N14").Value = 8.805885
$ws.Range("O14").Value = 0.0527739323334782
$ws.Range("P14").Value = 0.0527739323334782
$ws.Range("Q14").Value = 44.402744623985
$ws.Range("R14").Value = 399.624701615865
$ws.Range("S14").Value = 0.003157640342904481
$ws.Range("T14").Value = 0.003157640342904481

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.127183
$ws.Range("H15").Value = 45.381549
$ws.Range("I15").Value = 0.05983333443775553
$ws.Range("J15").Value = 0.05983333443775553
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 38.46708166666667
$ws.Range("N15").Value = 115.401245
$ws.Range("O15").Value = 0.69160311482936
$ws.Range("P15").Value = 0.69160311482936
$ws.Range("Q15").Value = 581.8985838476117
$ws.Range("R15").Value = 5237.087254628505
$ws.Range("S15").Value = 0.04138092046777853
$ws.Range("T15").Value = 0.04138092046777853

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.127183
$ws.Range("H16").Value = 45.381549
$ws.Range("I16").Value = 0.05983333443775553
$ws.Range("J16").Value = 0.05983333443775553
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.981753333333333
$ws.Range("N16").Value = 17.94526
$ws.Range("O16").Value = 0.1075464802171131
$ws.Range("P16").Value = 0.1075464802171131
$ws.Range("Q16").Value = 90.48707733419333
$ws.Range("R16").Value = 814.3836960077399
$ws.Range("S16").Value = 0.006434864518433985
$ws.Range("T16").Value = 0.006434864518433985

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.127183
$ws.Range("H17").Value = 45.381549
$ws.Range("I17").Value = 0.05983333443775553
$ws.Range("J17").Value = 0.05983333443775553
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.236038333333333
$ws.Range("N17").Value = 24.708115
$ws.Range("O17").Value = 0.1480764726200487
$ws.Range("P17").Value = 0.1480764726200487
$ws.Range("Q17").Value = 124.5880590633483
$ws.Range("R17").Value = 1121.292531570135
$ws.Range("S17").Value = 0.008859909108638523
$ws.Range("T17").Value = 0.008859909108638524
